# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values for rows 2-33, replacing the prior Strike# values in column G
$kValues = @(5, 3, 8, 8, 5, 4, 6, 5, 3, 6, 4, 8, 7, 8, 6, 4, 5, 5, 8, 8, 3, 8, 9, 7, 10, 9, 7, 7, 7, 5, 3, 3)

$row = 2
foreach ($val in $kValues) {
    $ws.Cells.Item($row, 7).Value = $val
    $row++
}
